{"js": "const body = context.document.body;\n\n// 1. Update the activation date from 2012 to 2023.\nconst dateResults = body.search(\"Ativa\u00e7\u00e3o: 01/01/2012\", { matchCase: true });\nawait context.sync();\ndateResults.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. \"Objetivos\" section: add an italic English translation paragraph right\n//    after the existing Portuguese paragraph.\nconst objetivos = body.search(\n  \"A disciplina visa propiciar aos alunos os conhecimentos b\u00e1sicos de eletroqu\u00edmica, tanto do ponto de vista da eletroqu\u00edmica i\u00f4nica como da eletr\u00f3dica, e apresentar as principais aplica\u00e7\u00f5es da eletroqu\u00edmica\",\n  { matchCase: true }\n);\nawait context.sync();\nconst objetivosTranslation = objetivos.items[0].insertParagraph(\n  \"The course aims to provide students with basic knowledge of electrochemistry, both from the point of view of ionic and electrodic electrochemistry, and to present the main applications of electrochemistry\",\n  Word.InsertLocation.after\n);\nobjetivosTranslation.font.italic = true;\nawait context.sync();\n\n// 3. \"Programa resumido\" section: add an italic English translation\n//    paragraph right after the existing Portuguese paragraph.\nconst resumido = body.search(\n  \"Princ\u00edpios da eletroqu\u00edmica i\u00f4nica e da eletroqu\u00edmica eletr\u00f3dica. Aplica\u00e7\u00f5es.\",\n  { matchCase: true }\n);\nawait context.sync();\nconst resumidoTranslation = resumido.items[0].insertParagraph(\n  \"Principles of ionic electrochemistry and electrodic electrochemistry. Applications.\",\n  Word.InsertLocation.after\n);\nresumidoTranslation.font.italic = true;\nawait context.sync();\n\n// 4. \"Programa\" section: add an italic English translation paragraph right\n//    after the existing Portuguese paragraph.\nconst programa = body.search(\n  \"Princ\u00edpios da eletroqu\u00edmica i\u00f4nica: intera\u00e7\u00f5es i\u00f4nicas, equil\u00edbrio i\u00f4nico e condu\u00e7\u00e3o eletrol\u00edtica. Princ\u00edpios da eletroqu\u00edmica eletr\u00f3dica: fen\u00f4menos interfaciais, potenciais de eletrodo e c\u00e9lulas eletroqu\u00edmicas. Processos de eletrodo. M\u00e9todos eletroqu\u00edmicos de an\u00e1lise qu\u00edmica. Aplica\u00e7\u00f5es da eletroqu\u00edmica: fontes eletroqu\u00edmicas de energia, processos eletrometal\u00fargicos e galvanoplastia.\",\n  { matchCase: true }\n);\nawait context.sync();\nconst programaTranslation = programa.items[0].insertParagraph(\n  \"Principles of ionic electrochemistry: ionic interactions, ionic equilibrium and electrolyte conduction. Principles of electrochemical electrochemistry: interfacial phenomena, electrode potentials and electrochemical cells. Electrode processes. Electrochemical methods of chemical analysis. Applications of electrochemistry: electrochemical sources of energy, electrometallurgical processes and electroplating.\",\n  Word.InsertLocation.after\n);\nprogramaTranslation.font.italic = true;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the activation date from 2012 to 2023.\n$d.Content.Find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2012\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ativa\u00e7\u00e3o: 01/01/2023\", 2) | Out-Null\n\n# 2. Insert italic English translations right after each Portuguese paragraph.\nfunction Insert-TranslationAfter($paraIndex, $text) {\n    $p = $d.Paragraphs($paraIndex)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs($paraIndex + 1)\n    $start = $newPara.Range.Start\n    $newPara.Range.InsertAfter($text)\n    $textRng = $d.Range($start, $start + $text.Length)\n    $textRng.Font.Italic = $true\n}\n\n# \"Objetivos\" translation (paragraph right after the Portuguese objectives text).\nInsert-TranslationAfter 6 \"The course aims to provide students with basic knowledge of electrochemistry, both from the point of view of ionic and electrodic electrochemistry, and to present the main applications of electrochemistry\"\n\n# \"Programa resumido\" translation.\nInsert-TranslationAfter 11 \"Principles of ionic electrochemistry and electrodic electrochemistry. Applications.\"\n\n# \"Programa\" translation.\nInsert-TranslationAfter 14 \"Principles of ionic electrochemistry: ionic interactions, ionic equilibrium and electrolyte conduction. Principles of electrochemical electrochemistry: interfacial phenomena, electrode potentials and electrochemical cells. Electrode processes. Electrochemical methods of chemical analysis. Applications of electrochemistry: electrochemical sources of energy, electrometallurgical processes and electroplating.\"\n"}
